# Natmi LR-pairs sheet update (Sdc2-Ptprj): rebuild the data rows to include the
# "ECs" sending-cluster group (Dr Hou advice) alongside the existing FAPs/sCs groups,
# recomputing every ligand/receptor expression + specificity metric per the new grouping.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ row=2; A="ECs"; B="Sdc2"; C="Ptprj"; D="ECs"; E=2; F=0.6666666666666666; G=1.321445333333333; H=3.964336; I=0.01021782062667047; J=0.01021782062667047; K=3; L=1; M=8.904209333333332; N=26.712628; O=0.3928890865119899; P=0.3928890865119899; Q=11.76642587055644; R=105.897832835008; S=0.004014470212155929; T=0.004014470212155929 },
  @{ row=3; A="ECs"; B="Sdc2"; C="Ptprj"; D="FAPs"; E=2; F=0.6666666666666666; G=1.321445333333333; H=3.964336; I=0.01021782062667047; J=0.01021782062667047; K=3; L=1; M=9.112038; N=27.336114; O=0.4020593128556135; P=0.4020593128556135; Q=12.041060092256; R=108.369540830304; S=0.004108169940041042; T=0.004108169940041043 },
  @{ row=4; A="ECs"; B="Sdc2"; C="Ptprj"; D="sCs"; E=2; F=0.6666666666666666; G=1.321445333333333; H=3.964336; I=0.01021782062667047; J=0.01021782062667047; K=3; L=1; M=4.64717; N=13.94151; O=0.2050516006323966; P=0.2050516006323966; Q=6.140981109706667; R=55.26882998736001; S=0.002095180474473497; T=0.002095180474473497 },
  @{ row=5; A="FAPs"; B="Sdc2"; C="Ptprj"; D="ECs"; E=3; F=1; G=105.9632263333333; H=317.889679; I=0.819340166699254; J=0.8193401666992541; K=3; L=1; M=8.904209333333332; N=26.712628; O=0.3928890865119899; P=0.3928890865119899; Q=943.5187489073791; R=8491.668740166411; S=0.3219098096370515; T=0.3219098096370515 },
  @{ row=6; A="FAPs"; B="Sdc2"; C="Ptprj"; D="FAPs"; E=3; F=1; G=105.9632263333333; H=317.889679; I=0.819340166699254; J=0.8193401666992541; K=3; L=1; M=9.112038; N=27.336114; O=0.4020593128556135; P=0.4020593128556135; Q=965.5409449519341; R=8689.868504567406; S=0.3294233444181059; T=0.3294233444181059 },
  @{ row=7; A="FAPs"; B="Sdc2"; C="Ptprj"; D="sCs"; E=3; F=1; G=105.9632263333333; H=317.889679; I=0.819340166699254; J=0.8193401666992541; K=3; L=1; M=4.64717; N=13.94151; O=0.2050516006323966; P=0.2050516006323966; Q=492.4291265194767; R=4431.86213867529; S=0.1680070126440966; T=0.1680070126440967 },
  @{ row=8; A="sCs"; B="Sdc2"; C="Ptprj"; D="ECs"; E=3; F=1; G=22.04284166666666; H=66.128525; I=0.1704420126740755; J=0.1704420126740755; K=3; L=1; M=8.904209333333332; N=26.712628; O=0.3928890865119899; P=0.3928890865119899; Q=196.2740765015222; R=1766.4666885137; S=0.06696480666278253; T=0.06696480666278254 },
  @{ row=9; A="sCs"; B="Sdc2"; C="Ptprj"; D="FAPs"; E=3; F=1; G=22.04284166666666; H=66.128525; I=0.1704420126740755; J=0.1704420126740755; K=3; L=1; M=9.112038; N=27.336114; O=0.4020593128556135; P=0.4020593128556135; Q=200.85521089465; R=1807.69689805185; S=0.06852779849746655; T=0.06852779849746657 },
  @{ row=10; A="sCs"; B="Sdc2"; C="Ptprj"; D="sCs"; E=3; F=1; G=22.04284166666666; H=66.128525; I=0.1704420126740755; J=0.1704420126740755; K=3; L=1; M=4.64717; N=13.94151; O=0.2050516006323966; P=0.2050516006323966; Q=102.4368325080833; R=921.93149257275; S=0.03494940751382639; T=0.0349494075138264 }
)

foreach ($r in $rows) {
  $ws.Range("A" + $r.row).Value = $r.A
  $ws.Range("B" + $r.row).Value = $r.B
  $ws.Range("C" + $r.row).Value = $r.C
  $ws.Range("D" + $r.row).Value = $r.D
  $ws.Range("E" + $r.row).Value = $r.E
  $ws.Range("F" + $r.row).Value = $r.F
  $ws.Range("G" + $r.row).Value = $r.G
  $ws.Range("H" + $r.row).Value = $r.H
  $ws.Range("I" + $r.row).Value = $r.I
  $ws.Range("J" + $r.row).Value = $r.J
  $ws.Range("K" + $r.row).Value = $r.K
  $ws.Range("L" + $r.row).Value = $r.L
  $ws.Range("M" + $r.row).Value = $r.M
  $ws.Range("N" + $r.row).Value = $r.N
  $ws.Range("O" + $r.row).Value = $r.O
  $ws.Range("P" + $r.row).Value = $r.P
  $ws.Range("Q" + $r.row).Value = $r.Q
  $ws.Range("R" + $r.row).Value = $r.R
  $ws.Range("S" + $r.row).Value = $r.S
  $ws.Range("T" + $r.row).Value = $r.T
}
